$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Permutation mapping: target row -> source row (source row supplies the "variable" field values)
$perm = @{}
$perm[2] = 97
$perm[3] = 58
$perm[4] = 83
$perm[5] = 66
$perm[6] = 29
$perm[7] = 100
$perm[8] = 64
$perm[9] = 92
$perm[10] = 45
$perm[11] = 46
$perm[12] = 6
$perm[13] = 7
$perm[14] = 51
$perm[15] = 101
$perm[16] = 109
$perm[17] = 85
$perm[18] = 75
$perm[19] = 72
$perm[20] = 90
$perm[21] = 22
$perm[22] = 23
$perm[23] = 73
$perm[24] = 54
$perm[25] = 86
$perm[26] = 40
$perm[27] = 44
$perm[28] = 95
$perm[29] = 52
$perm[30] = 39
$perm[31] = 35
$perm[32] = 36
$perm[33] = 19
$perm[34] = 114
$perm[35] = 53
$perm[36] = 69
$perm[37] = 60
$perm[38] = 88
$perm[39] = 32
$perm[40] = 33
$perm[41] = 37
$perm[42] = 34
$perm[43] = 13
$perm[44] = 24
$perm[45] = 74
$perm[46] = 26
$perm[47] = 49
$perm[48] = 50
$perm[49] = 93
$perm[50] = 94
$perm[51] = 113
$perm[52] = 81
$perm[53] = 104
$perm[54] = 77
$perm[55] = 56
$perm[56] = 78
$perm[57] = 68
$perm[58] = 102
$perm[59] = 111
$perm[60] = 21
$perm[61] = 79
$perm[62] = 57
$perm[63] = 105
$perm[64] = 43
$perm[65] = 47
$perm[66] = 48
$perm[67] = 20
$perm[68] = 61
$perm[69] = 62
$perm[70] = 28
$perm[71] = 112
$perm[72] = 18
$perm[73] = 5
$perm[74] = 2
$perm[75] = 31
$perm[76] = 59
$perm[77] = 80
$perm[78] = 9
$perm[79] = 82
$perm[80] = 110
$perm[81] = 65
$perm[82] = 41
$perm[83] = 42
$perm[84] = 67
$perm[85] = 12
$perm[86] = 106
$perm[87] = 63
$perm[88] = 103
$perm[89] = 55
$perm[90] = 91
$perm[91] = 87
$perm[92] = 3
$perm[93] = 71
$perm[94] = 108
$perm[95] = 14
$perm[96] = 10
$perm[97] = 11
$perm[98] = 27
$perm[99] = 17
$perm[100] = 4
$perm[101] = 89
$perm[102] = 98
$perm[103] = 30
$perm[104] = 25
$perm[105] = 76
$perm[106] = 84
$perm[107] = 107
$perm[108] = 70
$perm[109] = 38
$perm[110] = 99
$perm[111] = 15
$perm[112] = 8
$perm[113] = 16
$perm[114] = 96

# Columns that vary per row: D(4) Fecha, I(9) Calidad, J(10) Volumen, K(11) Precio minimo,
# L(12) Precio maximo, M(13) Precio promedio ponderado, N(14) Unidad de comercializacion,
# O(15) Origen, P(16) Precio $/Kg, Q(17) Kg o Unidades
$cols = @(4,9,10,11,12,13,14,15,16,17)

# Step 1: snapshot current values for all rows/cols into memory before overwriting anything
$orig = @{}
for ($r = 2; $r -le 114; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $orig[$r] = $rowVals
}

# Step 2: write back according to permutation (target row r receives data from source row perm[r])
for ($r = 2; $r -le 114; $r++) {
    $src = $perm[$r]
    $srcVals = $orig[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}

Write-Host "Reorder complete."